# TimeManage_ZYJ.xlsx - add two new log entries (2012.7.24 and 2012.7.25)
# "add some spc when storage del the ball and when fly speed is fast"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell {
    param($cell, [string]$text)
    # Prefix with an apostrophe so Excel doesn't auto-convert date-like
    # strings (e.g. "2012.7.24") into a date serial number, then strip
    # the resulting quote-prefix formatting so the cell keeps the sheet's
    # normal (General) style, just like the rest of the log rows.
    $cell.Value2 = "'" + $text
    $cell.ClearFormats()
}

# New row 54
Set-TextCell $ws.Range("A54") "2012.7.24"
Set-TextCell $ws.Range("B54") "加入每次得分显示，以及特效"
Set-TextCell $ws.Range("C54") "效果不是很满意，特别是特效"
$ws.Range("D54").Value = 4

# New row 55
Set-TextCell $ws.Range("A55") "2012.7.25"
Set-TextCell $ws.Range("B55") "加入削球特效。和飞行动物速度达到一定时特效效果一般"
$ws.Range("D55").Value = 3

# Match the author's final selection position recorded in the sheet view
$ws.Range("D58").Select() | Out-Null
